$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New submit row (5th submission) appended below the existing data.
$ws.Range("A6").Value = "5_201115_2117_rf_with_3in1_corrected_train_and_valid"
$ws.Range("B6").Value = 0.53
$ws.Range("C6").Value = "random foreset, in 3in1 data set, corrected train and valid set. Train (.84), valid (.82)"

# Column A widened to fit the longer "submit" label (author resized/auto-fit column A).
$ws.Columns.Item(1).ColumnWidth = 53.3

# Cursor ends up on the next empty row after entering the new comment.
$ws.Range("C7").Select() | Out-Null
